# Revisions based on reviewer comments.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)

    if ($shape.Name -eq "TextBox 3") {
        # "in-memory buffer" -> "In-memory buffer"
        $shape.TextFrame.TextRange.Text = "In-memory buffer"
    }
    elseif ($shape.Name -eq "TextBox 18") {
        # Fix typo: "Trasnport" -> "Transport"
        [void]$shape.TextFrame.TextRange.Replace("Trasnport", "Transport")
    }
    elseif ($shape.Name -eq "Rectangle 109") {
        # "message processing" -> "Message processing", with the box
        # re-centering/growing slightly (autofit) to match the new text width.
        $shape.TextFrame.TextRange.Text = "Message processing"
        $shape.Left = 496.295
        $shape.Width = 178.7689
    }
}
